$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 157.8
$ws.Range("I12").Value = 140
$ws.Range("J12").Value = 184.5
$ws.Range("K12").Value = 140
$ws.Range("L12").Value = 184.5
$ws.Range("M12").Value = 30
$ws.Range("N12").Value = -524.5
# Row 15
$ws.Range("H15").Value = 3325
$ws.Range("I15").Value = 3325
$ws.Range("K15").Value = 9975
$ws.Range("M15").Value = -9806
# Row 46
$ws.Range("H46").Value = 1000
$ws.Range("J46").Value = 1000
$ws.Range("L46").Value = 3000
$ws.Range("N46").Value = -3238
# Row 60
$ws.Range("H60").Value = 1000
$ws.Range("J60").Value = 1000
$ws.Range("L60").Value = 3000
$ws.Range("N60").Value = -3968
# Row 96
$ws.Range("H96").Value = 1119.2727
$ws.Range("I96").Value = 1376.6
$ws.Range("J96").Value = 567.8570999999999
$ws.Range("K96").Value = 4129.799999999999
$ws.Range("L96").Value = 1703.5713
$ws.Range("M96").Value = -2756.799999999999
$ws.Range("N96").Value = -4449.5713
# Row 127
$ws.Range("H127").Value = 922.3
$ws.Range("I127").Value = 689
$ws.Range("J127").Value = 1272.25
$ws.Range("K127").Value = 2067
$ws.Range("L127").Value = 3816.75
$ws.Range("M127").Value = 2893
$ws.Range("N127").Value = -13736.75
# Row 132
$ws.Range("H132").Value = 6804912.5
$ws.Range("I132").Value = 8132193
$ws.Range("J132").Value = 2598.25
$ws.Range("K132").Value = 24396579
$ws.Range("L132").Value = 7794.75
$ws.Range("M132").Value = -24394049
$ws.Range("N132").Value = -12854.75
# Row 135
$ws.Range("H135").Value = 133.63158
$ws.Range("I135").Value = 120.52941
$ws.Range("J135").Value = 245
$ws.Range("K135").Value = 1084.76469
$ws.Range("L135").Value = 2205
$ws.Range("M135").Value = 1450.23531
$ws.Range("N135").Value = -7275
# Row 137
$ws.Range("H137").Value = 1984.3
$ws.Range("I137").Value = 1133.6666
$ws.Range("K137").Value = 3400.9998
$ws.Range("M137").Value = -850.9998000000001
# Row 138
$ws.Range("H138").Value = 1717.59
$ws.Range("J138").Value = 1899.6941
$ws.Range("L138").Value = 5699.0823
$ws.Range("N138").Value = -15979.0823

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3450.58
$ws.Range("I32").Value = 3846.9
$ws.Range("J32").Value = 1865.3
$ws.Range("K32").Value = 3846.9
$ws.Range("L32").Value = 1865.3
$ws.Range("M32").Value = -3559.9
$ws.Range("N32").Value = -2439.3
# Row 61
$ws.Range("H61").Value = 864.8889
$ws.Range("I61").Value = 790.46155
$ws.Range("J61").Value = 2800
$ws.Range("K61").Value = 790.46155
$ws.Range("L61").Value = 2800
$ws.Range("M61").Value = -578.46155
$ws.Range("N61").Value = -3224
# Row 74
$ws.Range("H74").Value = 974.3929000000001
$ws.Range("I74").Value = 825.2273
$ws.Range("K74").Value = 825.2273
$ws.Range("M74").Value = 48.77269999999999
# Row 77
$ws.Range("H77").Value = 974.3929000000001
$ws.Range("I77").Value = 825.2273
$ws.Range("K77").Value = 4126.136500000001
$ws.Range("M77").Value = 241.8634999999995
# Row 110
$ws.Range("H110").Value = 1210.92
$ws.Range("I110").Value = 1098.0454
$ws.Range("J110").Value = 2038.6666
$ws.Range("K110").Value = 1098.0454
$ws.Range("L110").Value = 2038.6666
$ws.Range("M110").Value = 946.9546
$ws.Range("N110").Value = -6128.6666
# Row 136
$ws.Range("H136").Value = 864.8889
$ws.Range("I136").Value = 790.46155
$ws.Range("J136").Value = 2800
$ws.Range("K136").Value = 2371.38465
$ws.Range("L136").Value = 8400
$ws.Range("M136").Value = 178.61535
$ws.Range("N136").Value = -13500

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 9976.280000000001
$ws.Range("I134").Value = 7279.8945
$ws.Range("J134").Value = 18514.834
$ws.Range("K134").Value = 21839.6835
$ws.Range("L134").Value = 55544.50199999999
$ws.Range("M134").Value = -19304.6835
$ws.Range("N134").Value = -60614.50199999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1251.4474
$ws.Range("I31").Value = 833.52
$ws.Range("J31").Value = 2055.1538
$ws.Range("K31").Value = 833.52
$ws.Range("L31").Value = 2055.1538
$ws.Range("M31").Value = -538.52
$ws.Range("N31").Value = -2645.1538
# Row 34
$ws.Range("H34").Value = 1251.4474
$ws.Range("I34").Value = 833.52
$ws.Range("J34").Value = 2055.1538
$ws.Range("K34").Value = 833.52
$ws.Range("L34").Value = 2055.1538
$ws.Range("M34").Value = -631.52
$ws.Range("N34").Value = -2459.1538
# Row 99
$ws.Range("H99").Value = 1646549.9
$ws.Range("I99").Value = 2633304.8
$ws.Range("J99").Value = 1958.3334
$ws.Range("K99").Value = 2633304.8
$ws.Range("L99").Value = 1958.3334
$ws.Range("M99").Value = -2631806.8
$ws.Range("N99").Value = -4954.3334
# Row 126
$ws.Range("H126").Value = 1646549.9
$ws.Range("I126").Value = 2633304.8
$ws.Range("J126").Value = 1958.3334
$ws.Range("K126").Value = 7899914.399999999
$ws.Range("L126").Value = 5875.0002
$ws.Range("M126").Value = -7897444.399999999
$ws.Range("N126").Value = -10815.0002
# Row 132
$ws.Range("H132").Value = 8805.235000000001
$ws.Range("I132").Value = 10590.385
$ws.Range("J132").Value = 3003.5
$ws.Range("K132").Value = 31771.155
$ws.Range("L132").Value = 9010.5
$ws.Range("M132").Value = -29241.155
$ws.Range("N132").Value = -14070.5
# Row 134
$ws.Range("H134").Value = 11495553
$ws.Range("I134").Value = 13890022
$ws.Range("J134").Value = 2100
$ws.Range("K134").Value = 41670066
$ws.Range("L134").Value = 6300
$ws.Range("M134").Value = -41667531
$ws.Range("N134").Value = -11370

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1616.6666
$ws.Range("I5").Value = 2056.0715
$ws.Range("J5").Value = 737.8570999999999
$ws.Range("K5").Value = 6168.2145
$ws.Range("L5").Value = 2213.5713
$ws.Range("M5").Value = -6056.2145
$ws.Range("N5").Value = -2437.5713
# Row 68
$ws.Range("H68").Value = 2360.9443
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 2360.9443
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 7082.8329
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -8704.832900000001
# Row 71
$ws.Range("H71").Value = 2360.9443
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 2360.9443
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 21248.4987
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -29360.4987
# Row 82
$ws.Range("H82").Value = 9449.416999999999
# Row 85
$ws.Range("H85").Value = 9449.416999999999
# Row 98
$ws.Range("H98").Value = 2800
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 2800
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 8400
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -11396
# Row 131
$ws.Range("J131").Value = 1297.6769
$ws.Range("L131").Value = 3893.0307
$ws.Range("N131").Value = -13973.0307
# Row 135
$ws.Range("H135").Value = 1616.6666
$ws.Range("I135").Value = 2056.0715
$ws.Range("J135").Value = 737.8570999999999
$ws.Range("K135").Value = 18504.6435
$ws.Range("L135").Value = 6640.7139
$ws.Range("M135").Value = -15969.6435
$ws.Range("N135").Value = -11710.7139
# Row 141
$ws.Range("H141").Value = 1640.6666
$ws.Range("I141").Value = 1640.6666
$ws.Range("K141").Value = 4921.9998
$ws.Range("M141").Value = 258.0002000000004

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 123
$ws.Range("H123").Value = 10326
$ws.Range("J123").Value = 10326
$ws.Range("L123").Value = 10326
$ws.Range("N123").Value = -15226
# Row 126
$ws.Range("H126").Value = 3464.4119
$ws.Range("I126").Value = 2089.3333
$ws.Range("J126").Value = 4214.4546
$ws.Range("K126").Value = 6267.999899999999
$ws.Range("L126").Value = 12643.3638
$ws.Range("M126").Value = -3797.999899999999
$ws.Range("N126").Value = -17583.3638
# Row 128
$ws.Range("H128").Value = 40520
$ws.Range("J128").Value = 40520
$ws.Range("L128").Value = 40520
$ws.Range("N128").Value = -50480
# Row 132
$ws.Range("H132").Value = 1700.6857
$ws.Range("I132").Value = 1589.2142
$ws.Range("J132").Value = 2146.5715
$ws.Range("K132").Value = 4767.642599999999
$ws.Range("L132").Value = 6439.7145
$ws.Range("M132").Value = -2237.642599999999
$ws.Range("N132").Value = -11499.7145

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2138.75
$ws.Range("I7").Value = 2102.7778
$ws.Range("J7").Value = 2246.6667
$ws.Range("K7").Value = 2102.7778
$ws.Range("L7").Value = 2246.6667
$ws.Range("M7").Value = -1990.7778
$ws.Range("N7").Value = -2470.6667
# Row 59
$ws.Range("H59").Value = 15559.2
$ws.Range("J59").Value = 15559.2
$ws.Range("L59").Value = 15559.2
$ws.Range("N59").Value = -16867.2
# Row 126
$ws.Range("H126").Value = 2138.75
$ws.Range("I126").Value = 2102.7778
$ws.Range("J126").Value = 2246.6667
$ws.Range("K126").Value = 6308.3334
$ws.Range("L126").Value = 6740.000100000001
$ws.Range("M126").Value = -3838.3334
$ws.Range("N126").Value = -11680.0001
# Row 132
$ws.Range("H132").Value = 34821.1
$ws.Range("I132").Value = 1208.5186
$ws.Range("K132").Value = 3625.5558
$ws.Range("M132").Value = -1095.5558
# Row 136
$ws.Range("H136").Value = 7616.467
$ws.Range("I136").Value = 8074.7856
$ws.Range("K136").Value = 24224.3568
$ws.Range("M136").Value = -21674.3568
# Row 137
$ws.Range("H137").Value = 50414.5
$ws.Range("J137").Value = 50414.5
$ws.Range("L137").Value = 50414.5
$ws.Range("N137").Value = -60614.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 37038220
$ws.Range("I126").Value = 65360172
$ws.Range("K126").Value = 196080516
$ws.Range("M126").Value = -196078046
# Row 132
$ws.Range("H132").Value = 1580.9056
$ws.Range("I132").Value = 1378.5106
$ws.Range("K132").Value = 4135.531800000001
$ws.Range("M132").Value = -1605.531800000001
# Row 136
$ws.Range("H136").Value = 627.74286
$ws.Range("I136").Value = 326.05
$ws.Range("K136").Value = 978.1500000000001
$ws.Range("M136").Value = 1571.85
